$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the double space typo in the hospital name (B3) -> single space
$ws.Range("B3").Value = "부산대학교병원 진단검사의학과"

# Update the selected cell to B12 (next empty row after list), matching end-user interaction
$ws.Range("B12").Select()

# Trigger page setup so printer-settings relationship / pageSetup element gets written
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
